$wb = $excel.ActiveWorkbook

# --- Sheet1: remove row 3 (A3), which drops the dimension to A1:A2,
#     and move the selection to D17 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A3").EntireRow.Delete() | Out-Null
$ws1.Range("D17").Select() | Out-Null

# --- Sheet2: move the selection to A4 (single cell) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A4").Select() | Out-Null

# --- Add Sheet3 at the end of the workbook with the value that used
#     to live in Sheet1!A3 ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Sheet3"
$ws3.Range("A1").Value = "defect_trend.jpg"

$ws1.Select() | Out-Null
